$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B12").Value = "UnionOfTwoSortedArray"
$ws.Range("A12").Value = "Union and Intersection of two sorted arrays"

[void]$ws.Range("A12").Select()
